$d = $word.ActiveDocument

# 1) Title: "Versão 1.0" -> "Versão 2.0"
$d.Content.Find.Execute(
    "Versão 1.0", $true, $false, $false, $false, $false,
    $true, 1, $false, "Versão 2.0", 2) | Out-Null

# 2) Signature table (Table 1): row 2 is "Responsável pelo Projeto" / " Regina Albuquerque"
#    -> becomes "Patrocinadora" / "Regina Albuquerque" (leading space removed)
$sigTable = $d.Tables.Item(1)
$sigTable.Cell(2, 1).Range.Text = "Patrocinadora"
$sigTable.Cell(2, 2).Range.Text = "Regina Albuquerque"

# 3) Same table, row 4 ("Gabriel Martins") label changes from
#    "RQ (Representante da Qualidade)" -> "Responsável pelo Projeto"
$sigTable.Cell(4, 1).Range.Text = "Responsável pelo Projeto"

# 4) Table 2, row 2: the "Requisitos Doe em 5" storage-location cell gets the
#    URL text and loses its center alignment.
$locTable = $d.Tables.Item(2)
$locCell = $locTable.Cell(2, 2)
$locCell.Range.Text = "https://github.com/carloskrefer/RequisitosDoeEm5/blob/main/artefatos_avaliados/doe_5_requisitos.pdf "
$locCell.Range.Paragraphs.Item(1).Alignment = 0

# 5) Escalation sentence: "...para o Responsável pelo Projeto." -> "...para a Patrocinadora."
$d.Content.Find.Execute(
    "Caso a resolução das não conformidades não ocorra dentro do período definido na Seção 5, os avaliadores estarão escalonando elas para o Responsável pelo Projeto. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Caso a resolução das não conformidades não ocorra dentro do período definido na Seção 5, os avaliadores estarão escalonando elas para a Patrocinadora. ",
    2) | Out-Null
